$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Lis Indriani" record (row 6). Excel shifts the rows below it
# up by one, shrinking the tables (Table2/Table4/Table6/Table7) that spanned
# down to row 9 so they now end at row 8, and the used range shrinks from
# A1:N24 to A1:N23.
$ws.Rows(6).Delete()

# The payment total (previously 250000 across the four remaining students)
# needs to drop by the 50000 that belonged to the removed record.
$ws.Range("K9").Value = 200000

# Mirror the same corrected total in the recap table at the bottom of the
# sheet.
$ws.Range("F22").Value = 200000
$ws.Range("G22").Value = 200000
